$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B40").Value = 0.17
$ws.Range("C40").Value = 0.17
$ws.Range("E40").Value = 8
$ws.Range("G40").Value = "1x a cada 4 meses - irregular"
$ws.Range("H40").Value = 45817.73247685185
$ws.Range("I40").Value = 45939.73247685185
$ws.Range("J45").Value = "INATIVO - 1.9 meses sem comprar"
$ws.Range("C54").Value = 0.67
$ws.Range("D54").Value = 0.83
$ws.Range("E54").Value = 27
$ws.Range("F54").Value = 0.83
$ws.Range("H54").Value = 45817.5108912037
$ws.Range("I54").Value = 45878.5108912037
$ws.Range("E65").Value = 30
$ws.Range("H65").Value = 45817.94976851852
$ws.Range("I65").Value = 45847.94976851852
$ws.Range("J66").Value = "INATIVO - 11.9 meses sem comprar"
$ws.Range("J73").Value = "INATIVO - 32.8 meses sem comprar"
$ws.Range("J74").Value = "INATIVO - 7.1 meses sem comprar"
$ws.Range("J79").Value = "INATIVO - 21.6 meses sem comprar"
$ws.Range("J81").Value = "INATIVO - 25.6 meses sem comprar"
$ws.Range("J83").Value = "INATIVO - 21.1 meses sem comprar"
$ws.Range("J84").Value = "INATIVO - 8.8 meses sem comprar"
$ws.Range("J85").Value = "INATIVO - 14.4 meses sem comprar"
$ws.Range("J86").Value = "INATIVO - 4.3 meses sem comprar"
$ws.Range("J87").Value = "INATIVO - 11.1 meses sem comprar"
$ws.Range("J88").Value = "INATIVO - 10.5 meses sem comprar"
$ws.Range("J89").Value = "INATIVO - 14.4 meses sem comprar"
$ws.Range("J90").Value = "INATIVO - 32.8 meses sem comprar"
$ws.Range("J91").Value = "INATIVO - 13.1 meses sem comprar"
$ws.Range("J92").Value = "INATIVO - 18.1 meses sem comprar"
$ws.Range("J93").Value = "INATIVO - 15.8 meses sem comprar"
$ws.Range("J94").Value = "INATIVO - 18.5 meses sem comprar"
$ws.Range("J95").Value = "INATIVO - 32.3 meses sem comprar"
$ws.Range("J97").Value = "INATIVO - 1.6 meses sem comprar"
$ws.Range("J98").Value = "INATIVO - 22.2 meses sem comprar"
$ws.Range("J101").Value = "INATIVO - 13.9 meses sem comprar"
$ws.Range("E111").Value = 15399
$ws.Range("H111").Value = 45818.74002314815
$ws.Range("I111").Value = 45818.74002314815
